# Planilla de Métricas V2.1 Vagones de Tren - apply row 30 ("Ejecución de la Prueba")
# test-execution time entries that were previously left blank ("Completar").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Métricas")

# Fill in the test-execution timing data in row 30:
#   B30 = Tiempo Estimado (duration)   -> 00:15
#   C30 = Hora Inicio (time of day)    -> 17:40
#   D30 = Hora Fin (time of day)       -> 17:50
# E30 already holds a formula that computes D30-C30, it will recalc automatically.
# Assign plain numeric (date-serial) values so Excel keeps the existing cell
# number format/style instead of inferring a brand new one from a DateTime.
$ws.Range("B30").Value = 0.010416666666666666
$ws.Range("C30").Value = 0.73611111111111116
$ws.Range("D30").Value = 0.74305555555555547

# Recalculate the workbook so dependent formulas/cached chart values refresh.
$excel.Calculate()

# Restore the scroll position / selection that was active when the author saved the file.
$ws.Activate()
$appWindow = $excel.ActiveWindow
$appWindow.ScrollRow = 19
$appWindow.ScrollColumn = 1
$ws.Range("E30").Select()
